# Admitted.xlsx weekly update — adds the 04_05_2021 column (AA) with the
# new admission counts per age group, extends the "I alt" (total) row,
# and moves the visible window/selection to the newly added column,
# mirroring what a user would do by typing the new week's figures in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1): new date label for the added week -------------------
$ws.Range("AA1").Value = "04_05_2021"

# --- Data rows (2-11): admissions per age bracket for 04_05_2021 ---------
$ws.Range("AA2").Value  = 197
$ws.Range("AA3").Value  = 199
$ws.Range("AA4").Value  = 647
$ws.Range("AA5").Value  = 993
$ws.Range("AA6").Value  = 1440
$ws.Range("AA7").Value  = 2252
$ws.Range("AA8").Value  = 2283
$ws.Range("AA9").Value  = 3190
$ws.Range("AA10").Value = 2547
$ws.Range("AA11").Value = 702

# --- Totals row (12): sum of the new column, same pattern as B12:Z12 -----
$ws.Range("AA12").Formula = "=SUM(AA2:AA11)"

# --- View state: scroll right so the new column is in view and select it -
$win = $excel.ActiveWindow
$win.ScrollColumn = 18
$win.ScrollRow = 1
$ws.Range("Z16").Select()
